# Refined metadata to be additional tab
#
# 1. Refresh the "time_taken" timestamps in column F of the "data" sheet.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    PanelApp query that produced the data sheet.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Refresh data-sheet timestamps (column F, rows 2-40) -----------------
$newTimestamps = @(
    "2021-10-05 14:33:08.221705",
    "2021-10-05 14:33:08.221713",
    "2021-10-05 14:33:08.221717",
    "2021-10-05 14:33:08.221720",
    "2021-10-05 14:33:08.221723",
    "2021-10-05 14:33:08.221725",
    "2021-10-05 14:33:08.221728",
    "2021-10-05 14:33:08.221731",
    "2021-10-05 14:33:08.221734",
    "2021-10-05 14:33:08.221736",
    "2021-10-05 14:33:08.221739",
    "2021-10-05 14:33:08.221742",
    "2021-10-05 14:33:08.221744",
    "2021-10-05 14:33:08.221747",
    "2021-10-05 14:33:08.221750",
    "2021-10-05 14:33:08.221752",
    "2021-10-05 14:33:08.221755",
    "2021-10-05 14:33:08.221758",
    "2021-10-05 14:33:08.221760",
    "2021-10-05 14:33:08.221763",
    "2021-10-05 14:33:08.221766",
    "2021-10-05 14:33:08.221768",
    "2021-10-05 14:33:08.221771",
    "2021-10-05 14:33:08.221774",
    "2021-10-05 14:33:08.221777",
    "2021-10-05 14:33:08.221779",
    "2021-10-05 14:33:08.221782",
    "2021-10-05 14:33:08.221785",
    "2021-10-05 14:33:08.221787",
    "2021-10-05 14:33:08.221790",
    "2021-10-05 14:33:08.221792",
    "2021-10-05 14:33:08.221795",
    "2021-10-05 14:33:08.221798",
    "2021-10-05 14:33:08.221800",
    "2021-10-05 14:33:08.221803",
    "2021-10-05 14:33:08.221806",
    "2021-10-05 14:33:08.221808",
    "2021-10-05 14:33:08.221811",
    "2021-10-05 14:33:08.221813"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- 2. Add the "metadata" worksheet after "data" ----------------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1) - values then copy the bold/bordered header style used
# by the data sheet's header row (data!B1) so we reuse the existing style
# instead of minting a near-duplicate one.
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial($xlPasteFormats)

# Data row 2
$metaSheet.Cells.Item(2, 1).Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial($xlPasteFormats)

$metaSheet.Cells.Item(2, 2).Value = "Amelogenesis imperfecta"
$metaSheet.Cells.Item(2, 3).Value = 3564

# data_version must stay textual ("1.1") rather than become the number 1.1 -
# format the cell as Text before assigning, then restore the default
# (unstyled) look by pasting the data sheet's plain body format over it.
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.1"
$dataSheet.Range("B2").Copy()
$metaSheet.Range("D2").PasteSpecial($xlPasteFormats)

$metaSheet.Cells.Item(2, 5).Value = "2021-08-13T06:29:23.441723Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:33:08.217897"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3564/?format=json"

$excel.CutCopyMode = $false

# Keep "data" as the active sheet/selection, matching the workbook view
# before this edit (only the sheet list and cell content changed).
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A1").Select() | Out-Null

Write-Output "metadata sheet added; timestamps refreshed"
